# river update May 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Visual Clarity ---
$ws.Range("F2").Value = 0.138093282934764
$ws.Range("H2").Value = 0.827586206896552
$ws.Range("K2").Value = -0.0199046321525886
$ws.Range("L2").Value = -0.0492410937500995
$ws.Range("M2").Value = 0.0100026636955602
$ws.Range("N2").Value = -7.6556277509956
$ws.Range("P2").Value = "Unlikely improving"

# --- Row 3: Dissolved Oxygen Concentration ---
$ws.Range("F3").Value = 0.714906112454826
$ws.Range("H3").Value = 0.944444444444444
$ws.Range("J3").Value = 9.49
$ws.Range("K3").Value = 0.058575061029067
$ws.Range("L3").Value = -0.132776270173852
$ws.Range("M3").Value = 0.102789831370291
$ws.Range("N3").Value = 0.617229304837376
$ws.Range("P3").Value = "Likely increasing"

# --- Row 4: Dissolved Reactive Phosphorus ---
$ws.Range("F4").Value = 0.25555805922537
$ws.Range("H4").Value = 0.7037037037037041
$ws.Range("J4").Value = 0.0375
$ws.Range("K4").Value = 0.0011476826394344
$ws.Range("L4").Value = -0.0015333743789827
$ws.Range("M4").Value = 0.0039443553892329
$ws.Range("N4").Value = 3.06048703849175
$ws.Range("P4").Value = "Unlikely improving"

# --- Row 5: E. coli ---
$ws.Range("F5").Value = 0.938354201000995
$ws.Range("G5").Value = 0.0185185185185185
$ws.Range("H5").Value = 0.833333333333333
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = -35.9677796658383
$ws.Range("L5").Value = -56.2739561081485
$ws.Range("M5").Value = -4.17558282940239
$ws.Range("N5").Value = -14.3871118663353
$ws.Range("P5").Value = "Very likely improving"

# --- Row 6: Ammoniacal Nitrogen (NH4) ---
$ws.Range("E6").Value = "ok"
$ws.Range("F6").Value = 0.0346490789736528
$ws.Range("G6").Value = 0.351851851851852
$ws.Range("H6").Value = 0.7037037037037041
$ws.Range("J6").Value = 0.0175021742902809
$ws.Range("K6").Value = 0.0021980669161238
$ws.Range("M6").Value = 0.0055731148768861
$ws.Range("N6").Value = 12.5588220050147
$ws.Range("P6").Value = "Extremely unlikely improving"

# --- Row 7: Nitrite Nitrogen (NO2) ---
$ws.Range("D7").Value = $true
$ws.Range("F7").Value = 0.0347981289882045
$ws.Range("G7").Value = 0.12962962962963
$ws.Range("H7").Value = 0.592592592592593
$ws.Range("J7").Value = 0.0135
$ws.Range("K7").Value = 0.001126289536387
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0.0025038627113964
$ws.Range("N7").Value = 8.34288545471856
$ws.Range("P7").Value = "Extremely unlikely improving"

# --- Row 8: Nitrate Nitrogen (NO3) ---
$ws.Range("E8").Value = "WARNING: Sen slope influenced by censored values"
$ws.Range("F8").Value = 0.144243731661745
$ws.Range("G8").Value = 0.0925925925925926
$ws.Range("H8").Value = 0.907407407407407
$ws.Range("J8").Value = 0.353
$ws.Range("K8").Value = 0.0249506188613413
$ws.Range("L8").Value = -0.0297835040349404
$ws.Range("M8").Value = 0.0719089531463827
$ws.Range("N8").Value = 7.06816398338279
$ws.Range("P8").Value = "Unlikely improving"

# --- Row 9: pH ---
$ws.Range("F9").Value = 0.372395520847622
$ws.Range("H9").Value = 0.759259259259259
$ws.Range("J9").Value = 7.705
$ws.Range("K9").Value = -0.0200366980932163
$ws.Range("L9").Value = -0.07189080221147209
$ws.Range("M9").Value = 0.0300249850375088
$ws.Range("N9").Value = -0.260047996018382
$ws.Range("P9").Value = "As likely as not increasing"

# --- Row 10: SIN (Soluble Inorganic nitrogen) ---
$ws.Range("E10").Value = "WARNING: Sen slope influenced by censored values"
$ws.Range("F10").Value = 0.110335680959923
$ws.Range("G10").Value = 0.0925925925925926
$ws.Range("H10").Value = 0.685185185185185
$ws.Range("K10").Value = 0.0415475246890051
$ws.Range("L10").Value = -0.010918096922403
$ws.Range("M10").Value = 0.106230561284845
$ws.Range("N10").Value = 9.89226778309644

# --- Row 11: now Total Nitrogen (was Suspended Sediment Concentration) ---
$ws.Range("B11").Value = "Total Nitrogen"
$ws.Range("D11").Value = $false
$ws.Range("E11").Value = "ok"
$ws.Range("F11").Value = 0.002628382223813
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0.87037037037037
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 1.29
$ws.Range("K11").Value = 0.164527027027027
$ws.Range("L11").Value = 0.0722954880496792
$ws.Range("M11").Value = 0.254704022918368
$ws.Range("N11").Value = 12.7540331028703
$ws.Range("P11").Value = "Exceptionally unlikely improving"
$ws.Range("W11").Value = "g/m3"

# --- Row 12: now Total Phosphorus (was Total Nitrogen) ---
$ws.Range("B12").Value = "Total Phosphorus"
$ws.Range("F12").Value = 0.0017646835519445
$ws.Range("H12").Value = 0.907407407407407
$ws.Range("J12").Value = 0.138
$ws.Range("K12").Value = 0.0213123742454728
$ws.Range("L12").Value = 0.0100815224915003
$ws.Range("M12").Value = 0.0331159101459821
$ws.Range("N12").Value = 15.4437494532412
$ws.Range("P12").Value = "Exceptionally unlikely improving"

# --- Row 13: removed entirely ---
$ws.Rows.Item(13).Delete()
